# Hands On Demos - Day 4.
# Remove the stray "object 5" custom-geometry shape (a solid dark bar)
# from slide 14 of the "Nested Types and Anonymous Classes" deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$s.Shapes.Item("object 5").Delete()
